$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current layout (before edit):
#   row 179: 09-09-2024 | 943,18   (last data row)
#   row 180: (blank)
#   row 181: "Pie de página: Reporte generado automáticamente."
#   row 182: "Última actualización: 2024-09-08T07:38:50Z"
#
# Target layout (after edit):
#   row 179: 09-09-2024 | 943,18   (unchanged)
#   row 180: 10-09-2024 | 946,22   (new data row)
#   row 181: 11-09-2024 | 948,85   (new data row)
#   row 182: (blank, gap kept just like before)
#   row 183: "Pie de página: Reporte generado automáticamente."
#   row 184: "Última actualización: 2024-09-10T21:15:11Z"

# Move the footer text down first (into its final rows 183/184) so the
# rows it used to occupy (181/182) are free to be reused by new data.
$ws.Cells.Item(184, 1).Value = "Última actualización: 2024-09-10T21:15:11Z"
$ws.Cells.Item(183, 1).Value = "Pie de página: Reporte generado automáticamente."

# Clear the old footer cells before they get reused as data rows.
$ws.Cells.Item(181, 1).Value = $null
$ws.Cells.Item(182, 1).Value = $null

# The new dates look like dates ("10-09-2024"), so Excel would otherwise
# auto-convert them to date serial values on entry. Pre-format column A
# for these two rows as Text so they are stored verbatim as strings (like
# every other date cell in this sheet), then restore the Normal style so
# no formatting is left applied to the cells themselves.
$ws.Range("A180:A181").NumberFormat = "@"
$ws.Cells.Item(180, 1).Value = "10-09-2024"
$ws.Cells.Item(181, 1).Value = "11-09-2024"
$ws.Range("A180:A181").Style = "Normal"

$ws.Cells.Item(180, 2).Value = "946,22"
$ws.Cells.Item(181, 2).Value = "948,85"

$wb.Save()
